$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 84 (flight #83): Monday, Jan 16 - LO3804 to Warsaw
$ws.Range("A84").Value = 83
$ws.Range("B84").Value = "Monday, Jan 16"
$ws.Range("C84").Value = "5:40 AM"
$ws.Range("D84").Value = "LO3804"
$ws.Range("E84").Value = "Warsaw"
$ws.Range("F84").Value = "(WAW)"
$ws.Range("G84").Value = "LOT (Sliwka Naleczowska Livery) "
$ws.Range("H84").Value = "E195"
$ws.Range("I84").Value = "(SP-LNC)"
$ws.Range("J84").Value = "5:46 AM"
$ws.Range("L84").Value = "0 hours, 6 minutes"

# Row 85 (flight #84): Monday, Jan 16 - X7542 to Liege
$ws.Range("A85").Value = 84
$ws.Range("B85").Value = "Monday, Jan 16"
$ws.Range("C85").Value = "8:00 AM"
$ws.Range("D85").Value = "X7542"
$ws.Range("E85").Value = "Liege"
$ws.Range("F85").Value = "(LGG)"
$ws.Range("G85").Value = "Challenge Airlines "
$ws.Range("H85").Value = "B744"
$ws.Range("I85").Value = "(OO-ACE)"
$ws.Range("J85").Value = "8:26 AM"
$ws.Range("L85").Value = "0 hours, 26 minutes"
